# Refresh the cryptos list: updated prices/1h-deltas for existing rows,
# plus a new "WrappedliquidstakedEther2.0" entry inserted at row 33 which
# pushes every following coin down a row (dropping the old last row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.056.74'
$ws.Range('E2').Value = '  +2.14%  '
$ws.Range('D3').Value = '1.706.87'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '317.11'
$ws.Range('E5').Value = '  +0.28%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4002'
$ws.Range('E7').Value = '  +2.48%  '
$ws.Range('D8').Value = '0.4047'
$ws.Range('E8').Value = '  -0.90%  '
$ws.Range('D9').Value = '1.475'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').Value = '52.94'
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('D11').Value = '1.001'
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').Value = '0.08819'
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = '26.09'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').Value = '7.494'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '7.993'
$ws.Range('E15').Value = '  -4.02%  '
$ws.Range('D16').Value = '0.00001354'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '1.707.15'
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '96.09'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').Value = '0.07197'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.80'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '7.327'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = '14.33'
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('D24').Value = '25.043.21'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D26').Value = '2.946'
$ws.Range('E26').Value = '  -3.28%  '
$ws.Range('D27').Value = '23.58'
$ws.Range('E27').Value = '  +1.86%  '
$ws.Range('D28').Value = '6.079'
$ws.Range('E28').Value = '  +12.61%  '
$ws.Range('D29').Value = '162.94'
$ws.Range('E29').Value = '  -3.22%  '
$ws.Range('D30').Value = '152.37'
$ws.Range('E30').Value = '  +3.85%  '
$ws.Range('D31').Value = '8.429'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('D32').Value = '2.709'
$ws.Range('E32').Value = '  +23.36%  '
$ws.Range('B33').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C33').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D33').Value = '1.885.95'
$ws.Range('E33').Value = '  +0.27%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.08646'
$ws.Range('E34').Value = '  -2.04%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.03178'
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.049'
$ws.Range('E36').Value = '  -0.74%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '7.208'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2931'
$ws.Range('E38').Value = '  +4.21%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.09705'
$ws.Range('E39').Value = '  +5.54%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '11.02'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.8302'
$ws.Range('E41').Value = '  +3.57%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '14.06'
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.480'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '17.07'
$ws.Range('E44').Value = '  -2.40%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '2.694'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.7386'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.09223'
$ws.Range('E47').Value = '  +12.83%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.250'
$ws.Range('E48').Value = '  -0.45%  '
$ws.Range('B49').Value = 'Flow'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D49').Value = '1.409'
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').Value = '0.9999'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '140.12'
$ws.Range('E51').Value = '  -0.78%  '
